# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$targetSheets = @(1, 4)

foreach ($sheetIdx in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetIdx)

    $ws.Range("E2").Value = "2024.03.16 10:00-03.16 17:00"
    $ws.Range("E3").Value = "2024.03.16 09:30-03.17 17:00"
    $ws.Range("F3").Value = 3110
    $ws.Range("E4").Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Range("E5").Value = "2024.03.23 09:00-03.24 17:00"
    $ws.Range("F5").Value = 2785
    $ws.Range("E6").Value = "2024.03.23 10:00-03.24 17:00"
    $ws.Range("E7").Value = "2024.03.24 09:00-03.24 17:00"
    $ws.Range("F7").Value = 146
    $ws.Range("E8").Value = "2024.03.24 14:00-03.24 18:00"
    $ws.Range("E9").Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Range("F9").Value = 1493
    $ws.Range("E10").Value = "2024.03.30 10:00-03.30 17:00"
    $ws.Range("E11").Value = "2024.03.31 10:00-03.31 17:00"
    $ws.Range("E12").Value = "2024.03.31 14:00-03.31 18:00"
    $ws.Range("E13").Value = "2024.04.04 10:00-04.05 17:00"
    $ws.Range("F13").Value = 1244
    $ws.Range("E14").Value = "2024.04.04 10:00-04.04 17:00"
    $ws.Range("E15").Value = "2024.04.04 10:00-04.06 17:00"
    $ws.Range("F15").Value = 384
    $ws.Range("E16").Value = "2024.04.04 10:00-04.04 17:00"
    $ws.Range("E17").Value = "2024.04.05 09:30-04.05 17:00"
    $ws.Range("E18").Value = "2024.04.06 10:00-04.06 17:00"
    $ws.Range("E19").Value = "2024.04.13 10:00-04.13 17:00"
    $ws.Range("E20").Value = "2024.04.13 10:00-04.14 18:00"
    $ws.Range("F20").Value = 80
    $ws.Range("E21").Value = "2024.04.20 09:00-04.21 17:00"
    $ws.Range("E22").Value = "2024.04.20 09:30-04.21 17:00"
    $ws.Range("F22").Value = 2745
    $ws.Range("E23").Value = "2024.05.01 09:00-05.02 17:00"
    $ws.Range("E24").Value = "2024.05.02 10:00-05.02 16:00"
    $ws.Range("F24").Value = 6
    $ws.Range("E25").Value = "2024.05.26 09:30-05.26 17:30"
    $ws.Range("F25").Value = 47
}

Write-Output "edit complete"
